# Altera ordem processo SEI
#
# On the "compras-coronavirus" sheet, the "PROCESSO_SEI" / "URL_PROCESSO_SEI"
# entries (originally the first two data rows, 2-3) are moved to become the
# last two entries of the same mapping block (rows 21-22); every row that used
# to sit between them and the end of the block (old rows 4-22) shifts up by
# two positions. Column A (the running sequence number) and column D (a
# formula that always just mirrors column B on its own row) are left exactly
# as-is; D recalculates itself automatically once B changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compras-coronavirus")

$firstRow = 2
$lastRow = 22
$yellow = 65535

# 1) Snapshot the current B/C values (and whether C is highlighted yellow)
#    for every row in the affected block, before anything is overwritten.
$bVals = @{}
$cVals = @{}
$highlighted = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $bVals[$r] = $ws.Cells.Item($r, 2).Value()
    $cVals[$r] = $ws.Cells.Item($r, 3).Value()
    $highlighted[$r] = ($ws.Cells.Item($r, 3).Interior.Color -eq $yellow)
}

# 2) Work out, for every destination row, which row's old content now lands
#    there: old rows 4..22 move up to new rows 2..20, and old rows 2..3 wrap
#    around to become new rows 21..22.
$newOrder = @{}
$newRow = $firstRow
for ($r = $firstRow + 2; $r -le $lastRow; $r++) {
    $newOrder[$newRow] = $r
    $newRow++
}
for ($r = $firstRow; $r -le $firstRow + 1; $r++) {
    $newOrder[$newRow] = $r
    $newRow++
}

# 3) Write the rotated values back, carrying the yellow highlight along with
#    whichever label now occupies each row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $newOrder[$r]
    $ws.Cells.Item($r, 2).Value = $bVals[$src]
    $ws.Cells.Item($r, 3).Value = $cVals[$src]

    if ($highlighted[$src]) {
        $ws.Cells.Item($r, 3).Interior.Color = $yellow
    } else {
        $ws.Cells.Item($r, 3).Style = "Normal"
    }
}

# 4) Refresh the remembered selection/scroll state on "compras-coronavirus"
#    (it must stay a non-active tab, matching the saved workbook).
$ws.Activate()
$ws.Range("B32:B33").Select()

# 5) Refresh the remembered selection/zoom on "mapeamento", and leave it as
#    the active tab, matching the saved workbook.
$mapSheet = $wb.Worksheets.Item("mapeamento")
$mapSheet.Activate()
$excel.ActiveWindow.Zoom = 110
$mapSheet.Range("C1").Select()
$mapSheet.Range("C2").Select()
